# ----------------------------------------------------------------------
# Add a new "Seat Assignments" sheet (position 3, right after
# "Contestants"), renumber the sheets that follow it, seed the new
# sheet with its header row + one seat-assignment record, and update
# the Contestants sheet so Felicity Parker-Hill / Peter Adamidis swap
# rows (Felicity -> row 2, now "assigned"; Peter -> row 4).
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Snapshot the data of the four sheets that sit after "Contestants"
#    so we can recreate them (in order) once the new sheet has been
#    inserted. Re-creating them (instead of just moving them) is what
#    gives the new "Seat Assignments" sheet sheetId=3 and pushes the
#    others to 4/5/6/7, matching a freshly authored workbook.
# ------------------------------------------------------------------
$standbysVals     = $wb.Worksheets.Item("Standbys").Range("A1:E2").Value2
$groupsVals       = $wb.Worksheets.Item("Groups").Range("A1:B6").Value2
$blockTypesVals   = $wb.Worksheets.Item("Block Types").Range("A1:D27").Value2
$canceledVals     = $wb.Worksheets.Item("Canceled Assignments").Range("A1:D2").Value2
$canceledDate     = $wb.Worksheets.Item("Canceled Assignments").Range("E2").Value2

# ------------------------------------------------------------------
# 2. Drop the four sheets - this frees up their sheetIds (3-6) so the
#    next sheet we add (Seat Assignments) claims sheetId 3.
# ------------------------------------------------------------------
$wb.Worksheets.Item("Canceled Assignments").Delete()
$wb.Worksheets.Item("Block Types").Delete()
$wb.Worksheets.Item("Groups").Delete()
$wb.Worksheets.Item("Standbys").Delete()

# ------------------------------------------------------------------
# 3. Recreate each sheet (in final left-to-right order), moving it to
#    the end of the tab strip right after it is added so the order
#    comes out: Record Days, Contestants, Seat Assignments, Standbys,
#    Groups, Block Types, Canceled Assignments.
# ------------------------------------------------------------------

# NOTE: each sheet is fully populated with its data *before* the
# `.Move()` call that relocates it to the end of the tab strip - doing
# the write after the move silently fails to persist on save.

# -- Seat Assignments (brand-new sheet + data) --
$seat = $wb.Worksheets.Add()
$seat.Name = "Seat Assignments"

$seat.Range("A1").Value = "ID"
$seat.Range("B1").Value = "RecordDayID"
$seat.Range("C1").Value = "ContestantID"
$seat.Range("D1").Value = "Block"
$seat.Range("E1").Value = "Seat"
$seat.Range("F1").Value = "BookingEmailSent"
$seat.Range("G1").Value = "ConfirmedRSVP"
$seat.Range("H1").Value = "Notes"

$seat.Range("A2").Value = "ccd6460c-a021-4806-8e21-22177e70d8ba"
$seat.Range("B2").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$seat.Range("C2").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$seat.Range("D2").Value = 1
$seat.Range("E2").Value = "B3"

$seat.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# -- Standbys (restored) --
$standbys = $wb.Worksheets.Add()
$standbys.Name = "Standbys"
$standbys.Range("A1:E2").Value = $standbysVals
$standbys.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# -- Groups (restored) --
$groups = $wb.Worksheets.Add()
$groups.Name = "Groups"
$groups.Range("A1:B6").Value = $groupsVals
$groups.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# -- Block Types (restored) --
$blockTypes = $wb.Worksheets.Add()
$blockTypes.Name = "Block Types"
$blockTypes.Range("A1:D27").Value = $blockTypesVals
$blockTypes.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# -- Canceled Assignments (restored, incl. the formatted date cell) --
$canceled = $wb.Worksheets.Add()
$canceled.Name = "Canceled Assignments"
$canceled.Range("A1:D2").Value = $canceledVals
$canceled.Range("E1").Value = "CanceledAt"
$canceled.Range("E2").Value = $canceledDate
$canceled.Range("E2").NumberFormat = "m/d/yyyy"
$canceled.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# ------------------------------------------------------------------
# 4. Contestants sheet: Felicity Parker-Hill and Peter Adamidis swap
#    rows. Felicity moves to row 2 and becomes "assigned" (she now has
#    a seat); Peter moves to row 4, keeping his "assigned" status.
# ------------------------------------------------------------------
$contestants = $wb.Worksheets.Item("Contestants")

$contestants.Range("A2").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$contestants.Range("B2").Value = "Felicity Parker-Hill"
$contestants.Range("C2").Value = 27
$contestants.Range("D2").Value = "Not Specified"
$contestants.Range("E2").Value = "felicity.parkerhill@endemolshine.com.au"
$contestants.Range("F2").Value = "498086080"
$contestants.Range("G2").Value = "Melbourne"
$contestants.Range("H2").Value = ""
$contestants.Range("I2").Value = "assigned"
$contestants.Range("J2").Value = "Peter Adamidis, Kathleen Reynolds"
$contestants.Range("K2").Value = "5fe641da-4067-49a7-bae7-e63413b3e404"
$contestants.Range("L2").Value = "N"
$contestants.Range("M2").Value = "N/A"

$contestants.Range("A4").Value = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"
$contestants.Range("B4").Value = "Peter Adamidis"
$contestants.Range("C4").Value = 34
$contestants.Range("D4").Value = "Not Specified"
$contestants.Range("E4").Value = "peter.adamidis@gmail.com"
$contestants.Range("F4").Value = "498086080"
$contestants.Range("G4").Value = ""
$contestants.Range("H4").Value = ""
$contestants.Range("I4").Value = "assigned"
$contestants.Range("J4").Value = "Kathleen Reynolds, Felicity Parker-Hill"
$contestants.Range("K4").Value = "5fe641da-4067-49a7-bae7-e63413b3e404"
$contestants.Range("L4").Value = "Y"
$contestants.Range("M4").Value = "Broken Leg"

Write-Output "done"
